$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New futures rows for "Week 5", appended after the existing data (rows 2-201).
$newRows = @(
    @(340, "Ohio State"),
    @(360, "Georgia"),
    @(500, "Texas"),
    @(800, "Alabama"),
    @(1000, "Oregon"),
    @(1200, "Tennessee"),
    @(1200, "Ole Miss"),
    @(1800, "Miami"),
    @(2000, "Penn State"),
    @(4000, "Clemson"),
    @(5000, "LSU"),
    @(5000, "Utah"),
    @(5000, "USC"),
    @(5000, "Notre Dame"),
    @(5000, "Missouri"),
    @(8000, "Texas A&M"),
    @(12500, "Kansas State"),
    @(12500, "Michigan"),
    @(12500, "Oklahoma"),
    @(12500, "Louisville"),
    @(20000, "Oklahoma State"),
    @(25000, "Boise State"),
    @(25000, "Iowa State"),
    @(30000, "Nebraska"),
    @(35000, "Auburn"),
    @(30000, "UCF"),
    @(40000, "Iowa"),
    @(40000, "Indiana"),
    @(40000, "Kentucky"),
    @(50000, "Illinois"),
    @(50000, "Washington"),
    @(50000, "Washington State"),
    @(100000, "Arizona State"),
    @(100000, "Liberty"),
    @(10000, "Rutgers"),
    @(50000, "Maryland"),
    @(100000, "UNLV"),
    @(40000, "TCU"),
    @(100000, "Fresno State"),
    @(50000, "Georgia Tech"),
    @(40000, "Virginia Tech"),
    @(100000, "California"),
    @(100000, "Texas Tech"),
    @(150000, "Texas State"),
    @(100000, "Wisconsin"),
    @(50000, "South Carolina"),
    @(50000, "West Virginia"),
    @(50000, "Arizona"),
    @(100000, "Pittsburgh"),
    @(50000, "Colorado")
)

$startRow = 202
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $odds = $newRows[$i][0]
    $team = $newRows[$i][1]
    $ws.Cells.Item($r, 1).Value = $odds
    $ws.Cells.Item($r, 2).Value = $team
    $ws.Cells.Item($r, 3).Value = 5
}

$ws.Range("D234").Select()

"Added $($newRows.Count) rows"
